$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

$bfMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a6afa43e2d02990010dff88c76c3043d6a140c2/e2e/bf6bf436-adc4-49c7-9d72-96a56a3ec95a.md"
$bfMdName = "bf6bf436-adc4-49c7-9d72-96a56a3ec95a.md"

# --- 2. zh-cn sheet: fill in "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $bfMdUrl, "", "", $bfMdName)
$wsZh.Range("J2").Value2 = $wsZh.Range("G2").Value2
$wsZh.Range("K2").Value = "2016-09-07 13:52:19"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $bfMdUrl, "", "", $bfMdName)
$wsZh.Range("J3").Value2 = $wsZh.Range("G3").Value2
$wsZh.Range("K3").Value = "2016-09-07 13:52:19"

# --- 3. de-de sheet: same but with its own handback datetime ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $bfMdUrl, "", "", $bfMdName)
$wsDe.Range("J2").Value2 = $wsDe.Range("G2").Value2
$wsDe.Range("K2").Value = "2016-09-07 13:52:41"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $bfMdUrl, "", "", $bfMdName)
$wsDe.Range("J3").Value2 = $wsDe.Range("G3").Value2
$wsDe.Range("K3").Value = "2016-09-07 13:52:41"

# --- 4. Column width adjustments (report now shows longer status text / longer hyperlinks) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 29.9777050018311
$wsOverview.Range("F1").ColumnWidth = 29.9777050018311

foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Range("C1").ColumnWidth = 29.9777050018311
    $ws.Range("I1").ColumnWidth = 40
    $ws.Range("J1").ColumnWidth = 40
}

Write-Host "Generated handback report"
